$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Fill in remaining cells of existing row 10 (X10, Y10) ---
$ws.Range("X10").Value = -0.010002000000000066
$ws.Range("Y10").Value = "Down"

# --- Append new row 11 with a full set of trading data ---
$ws.Range("A11").Value = 42654.882118055553
$ws.Range("A11").NumberFormat = "m/d/yy h:mm"

$ws.Range("B11").Value = 3
$ws.Range("C11").Value = "Neutral"
$ws.Range("D11").Value = 14
$ws.Range("E11").Value = 30004
$ws.Range("F11").Value = 3419
$ws.Range("G11").Value = 59
$ws.Range("H11").Value = 35
$ws.Range("I11").Value = 80
$ws.Range("J11").Value = 19
$ws.Range("K11").Value = 16210
$ws.Range("L11").Value = 394
$ws.Range("M11").Value = 239
$ws.Range("N11").Value = 96
$ws.Range("O11").Value = 23
$ws.Range("P11").Value = "Bag"
$ws.Range("Q11").Value = 17.089518681678967
$ws.Range("R11").Value = -24.44

$ws.Range("S11").Value = -0.1101
$ws.Range("S11").NumberFormat = $ws.Range("S10").NumberFormat

$ws.Range("T11").Value = -0.0419
$ws.Range("T11").NumberFormat = $ws.Range("T10").NumberFormat

$ws.Range("U11").Value = 6.47
$ws.Range("V11").Value = 1.88
$ws.Range("W11").Value = -2
